$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F58").Value = '매 학년도 전기 졸업식은 2월 중순, 후기 졸업식은 8월 중순에 진행. 정확한 일정은 학교 홈페이지 학사공지 https://www.hknu.ac.kr/kor/562/subview.do 참고 또는 학사지원팀 031-670-5036 연락'
$ws.Range("F50").Value = '원전공은 전필과목 포함 39학점 이상 이수하여야 하고 복수전공은 전필과목 포함 39학점 이상 이수, 법학전공을 제외한 인문사회계열은 36학점 이상 이수. 교양은 원전공의 교양 졸업요구학점 이상 이수하면 되고, 졸업인증은 원전공만 합격, 졸업종합시험은 원전공 및 복수전공 둘다 합격해야 함. 인문계열 전공 학생이 자연/공학계열 전공을 선택하면 총 졸업학점 수가 120학점에서 130학점으로 변경됨.'
$ws.Range("F51").Value = '원전공은 전필과목 포함 51학점 이상 이수하여야 하고 부전공은 전필과목 포함 27학점 이상 이수, 원전공은 법학전공을 제외한 인문사회계열은 48학점 이상 이수. 교양은 원전공의 교양 졸업요구학점 이상 이수하면 되고, 졸업인증, 졸업종합시험은 원전공만 합격하면 됨.'
$ws.Range("F52").Value = '원전공은 전필과목 포함 39학점 이상 이수하여야 하고 융합전공은 전필과목 포함 39학점 이상 이수, 원전공은 법학전공을 제외한 인문사회계열은 36학점 이상 이수. 교양은 원전공의 교양 졸업요구학점 이상 이수하면 되고, 졸업인증은 원전공만 합격, 졸업종합시험은 원전공 및 융합전공 둘다 합격해야 함.'
$ws.Range("F53").Value = '원전공은 전필과목 포함 51학점 이상 이수하여야 하고 융합부전공은 전필과목 포함 27학점 이상 이수, 원전공은 법학전공을 제외한 인문사회계열은 48학점 이상 이수. 교양은 원전공의 교양 졸업요구학점 이상 이수하면 되고, 졸업인증, 졸업종합시험은 원전공만 합격하면 됨.'
$ws.Range("F54").Value = '원전공은 전필과목 포함 78학점 이상 이수하여야 하고 연계전공은 전공 교과목 9학점 이상 이수, 원전공은 법학전공을 제외한 인문사회계열은 72학점 이상 이수. 교양은 원전공의 교양 졸업요구학점 이상 이수하면 되고, 졸업인증, 졸업종합시험은 원전공만 합격하면 됨.'
$ws.Range("F49").Value = '1. 복수전공 제도는 원전공 및 복수전공의 전필과목 포함 39학점 이상 이수. 법학전공을 제외한 인문사회계열은 36학점 이상 이수, 2. 부전공 제도는 원전공은 전필과목 포함 51학점 이상 이수, 법학전공을 제외한 인문사회계열은 48학점 이상 이수. 부전공은 전필과목 포함 27학점 이상 이수, 3. 융합전공 제도는 원전공은 전필과목 포함 39학점 이상 이수, 법학전공을 제외한 인문사회계열은 36학점 이상 이수. 융합전공은 전필과목 포함 39학점 이상 이수, 4. 융합부전공 제도는 원전공은 전필과목 포함 51학점 이상 이수, 법학전공을 제외한 인문사회계열은 48학점 이상 이수. 융합부전공 제도는 전필과목 포함 27학점 이상 이수, 5. 연계전공 제도는 연계전공의 교과목 9학점 이상 이수, 6. 소단위전공과정(마이크로디그리) 제도는 과정별로 정해진 3~4과목 이수. 전공 신청 전에 이수한 교과목도 해당 전공으로 소급 인정됨. 공통으로 성적증명서 상단에 표시가 되며, 학위증에 다전공명 표시됨. 전공별 세부사항은 왼쪽 메뉴의 다전공 제도 안내 참고.'
$ws.Range("D42").Value = '다전공제도, 다전공뭐, 다전공제도뭐, 다전공무엇, 다전공제도무엇, 다전공알려, 다전공제도알려, 다전공설명, 다전공제도에대해설명, 다전공제도에대해알려, 다전공의개념, 다전공개념, 다전공정의, 다전공의정의가뭐, 다전공이뭐, 유연학사제도가뭐, 유연학사제도뭐, 유연학사제도에대해설명, 유연학사제도가무엇, 유연학사제도가뭐, 유연학사제도개념, 유연학사제도의개념, 유연학사제도의정의, 유연학사제도정의, 유연학사제도알려, 유연학사제도에대해알려, 유연학사제도설명'
$ws.Range("D56").Value = '마이크로디그리학점, 마이크로디그리몇학점, 마이크로디그리이수, 마이크로디그리이수학점, 마이크로디그리필수학점, 마이크로디그리졸업학점, 마이크로디그리는몇학점, 마이크로디그리의이수학점, 마이크로디그리과정의이수학점'

$ws.Application.Goto($ws.Range("A35"), $true)
$ws.Range("D56").Select() | Out-Null
